$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert 9 new rows at row 11 (pushes the existing "US..Argentina" blocks
#    down by 9 rows and auto-adjusts the C2:C10 SUM() formulas + dimension).
$ws.Rows("11:19").Insert()

# 2. Populate the new rows 11-19 with the "World" block.
$worldCounts = @(5928, 2707, 1767, 1337, 1019, 854, 714, 639, 525)
$worldExpected = @("301/1000", "176/1000", "125/1000", "97/1000", "79/1000", "67/1000", "58/1000", "46/1000")
$expectedFormulas = @("301/1000", "176/1000", "125/1000", "97/1000", "79/1000", "67/1000", "58/1000", "51/1000", "46/1000")

for ($i = 0; $i -lt 9; $i++) {
    $r = 11 + $i
    $ws.Cells.Item($r, 1).Value2 = "World"
    $ws.Cells.Item($r, 2).Value2 = $i + 1
    $ws.Cells.Item($r, 3).Value2 = $worldCounts[$i]
    $ws.Cells.Item($r, 4).Formula = "=SUM(C11:C19)"
    $ws.Cells.Item($r, 5).Formula = "=(C" + $r + "/D" + $r + ")"
    $ws.Cells.Item($r, 6).Formula = "=" + $expectedFormulas[$i]
    $ws.Cells.Item($r, 7).Formula = "=(E" + $r + "-F" + $r + ")"
}

# 3. Fix up the summary rows (2-10): label + total formulas.
$ws.Range("A2:A10").Value2 = "15_summed"
$ws.Cells.Item(1, 4).Value2 = "World"

$countryRowsOld = @(11,20,29,38,47,56,65,74,83,92,101,110,119,128,137)
for ($r = 2; $r -le 10; $r++) {
    $terms = @()
    foreach ($base in $countryRowsOld) {
        $terms += ("C" + ($base + $r - 2))
    }
    $ws.Cells.Item($r, 3).Formula = "=SUM(" + ($terms -join ",") + ")"
}

# 4. Trim the conditional-formatting dxfs back down to 2 entries (drop the
#    3 duplicate entries that used to back the old extra CF rules) and
#    repoint/resize the surviving rule to the new used range.
$ws.Cells.Item(1,1).Select() | Out-Null
